# Rename existing sheets
$wb = $excel.ActiveWorkbook
$wsJune = $wb.Worksheets.Item(1)
$wsSept = $wb.Worksheets.Item(2)
$wsJune.Name = "June_22 "
$wsSept.Name = "Sept_22"

# Add the new quarterly sheets after Sept_22, in order
$wsJan23  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSept)
$wsJan23.Name = "Jan_23"
$wsMar23  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsJan23)
$wsMar23.Name = "Mar_23"
$wsJune23 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsMar23)
$wsJune23.Name = "June_23"
$wsSept23 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsJune23)
$wsSept23.Name = "Sept_23"
$wsJan24  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSept23)
$wsJan24.Name = "Jan_24"

# ---- June_22 sheet: column widths + selection ----
$wsJune.Columns.Item(4).ColumnWidth = 15.0
$wsJune.Columns.Item(5).ColumnWidth = 13.666666666666666
$wsJune.Columns.Item(6).ColumnWidth = 17.5
$wsJune.Columns.Item(7).ColumnWidth = 19.0

# ---- Sept_22 sheet: column width ----
$wsSept.Columns.Item(6).ColumnWidth = 13.833333333333334

# ---- Jan_23 sheet: data ----
$jan23Data = New-Object 'object[,]' 63,6
$jan23Data[0,0] = 'District'
$jan23Data[0,1] = 'Sector'
$jan23Data[0,2] = 'totalDiagnoses'
$jan23Data[0,3] = 'totalIncidence'
$jan23Data[0,4] = 'incidenceRate'
$jan23Data[0,5] = 'Level'
$jan23Data[1,0] = 'RULINDO'
$jan23Data[1,1] = 'GISHYITA'
$jan23Data[1,2] = 25
$jan23Data[1,3] = 0
$jan23Data[1,4] = 0
$jan23Data[1,5] = 'Low'
$jan23Data[2,0] = 'KAYONZA'
$jan23Data[2,1] = 'NDEGO'
$jan23Data[2,2] = 36
$jan23Data[2,3] = 0
$jan23Data[2,4] = 0
$jan23Data[2,5] = 'Low'
$jan23Data[3,0] = 'RUBAVU'
$jan23Data[3,1] = 'CYANZARWE'
$jan23Data[3,2] = 80
$jan23Data[3,3] = 0
$jan23Data[3,4] = 0
$jan23Data[3,5] = 'Low'
$jan23Data[4,0] = 'KARONGI'
$jan23Data[4,1] = 'MURAMBI'
$jan23Data[4,2] = 123
$jan23Data[4,3] = 0
$jan23Data[4,4] = 0
$jan23Data[4,5] = 'Low'
$jan23Data[5,0] = 'KARONGI'
$jan23Data[5,1] = 'GISHYITA'
$jan23Data[5,2] = 307
$jan23Data[5,3] = 2
$jan23Data[5,4] = 0.65
$jan23Data[5,5] = 'Low'
$jan23Data[6,0] = 'BURERA'
$jan23Data[6,1] = 'KINYABABA'
$jan23Data[6,2] = 141
$jan23Data[6,3] = 1
$jan23Data[6,4] = 0.71
$jan23Data[6,5] = 'Low'
$jan23Data[7,0] = 'KAYONZA'
$jan23Data[7,1] = 'RUKARA'
$jan23Data[7,2] = 562
$jan23Data[7,3] = 4
$jan23Data[7,4] = 0.71
$jan23Data[7,5] = 'Low'
$jan23Data[8,0] = 'GISAGARA'
$jan23Data[8,1] = 'MUGANZA'
$jan23Data[8,2] = 54
$jan23Data[8,3] = 1
$jan23Data[8,4] = 1.85
$jan23Data[8,5] = 'Low'
$jan23Data[9,0] = 'GATSIBO'
$jan23Data[9,1] = 'REMERA'
$jan23Data[9,2] = 100
$jan23Data[9,3] = 2
$jan23Data[9,4] = 2
$jan23Data[9,5] = 'Low'
$jan23Data[10,0] = 'BURERA'
$jan23Data[10,1] = 'RUGARAMA'
$jan23Data[10,2] = 235
$jan23Data[10,3] = 6
$jan23Data[10,4] = 2.5499999999999998
$jan23Data[10,5] = 'Low'
$jan23Data[11,0] = 'BURERA'
$jan23Data[11,1] = 'KINONI'
$jan23Data[11,2] = 29
$jan23Data[11,3] = 1
$jan23Data[11,4] = 3.45
$jan23Data[11,5] = 'Low'
$jan23Data[12,0] = 'HUYE'
$jan23Data[12,1] = 'HUYE'
$jan23Data[12,2] = 28
$jan23Data[12,3] = 1
$jan23Data[12,4] = 3.57
$jan23Data[12,5] = 'Low'
$jan23Data[13,0] = 'BURERA'
$jan23Data[13,1] = 'NDAGO'
$jan23Data[13,2] = 28
$jan23Data[13,3] = 1
$jan23Data[13,4] = 3.57
$jan23Data[13,5] = 'Low'
$jan23Data[14,0] = 'RUTSIRO'
$jan23Data[14,1] = 'MUSASA'
$jan23Data[14,2] = 23
$jan23Data[14,3] = 1
$jan23Data[14,4] = 4.3499999999999996
$jan23Data[14,5] = 'Low'
$jan23Data[15,0] = 'GICUMBI'
$jan23Data[15,1] = 'MUTETE'
$jan23Data[15,2] = 23
$jan23Data[15,3] = 1
$jan23Data[15,4] = 4.3499999999999996
$jan23Data[15,5] = 'Low'
$jan23Data[16,0] = 'RUTSIRO'
$jan23Data[16,1] = 'BONEZA'
$jan23Data[16,2] = 64
$jan23Data[16,3] = 3
$jan23Data[16,4] = 4.6900000000000004
$jan23Data[16,5] = 'Low'
$jan23Data[17,0] = 'MUHANGA'
$jan23Data[17,1] = 'NYAMABUYE'
$jan23Data[17,2] = 86
$jan23Data[17,3] = 5
$jan23Data[17,4] = 5.81
$jan23Data[17,5] = 'Low'
$jan23Data[18,0] = 'BURERA'
$jan23Data[18,1] = 'RUSARABUYE'
$jan23Data[18,2] = 621
$jan23Data[18,3] = 37
$jan23Data[18,4] = 5.96
$jan23Data[18,5] = 'Low'
$jan23Data[19,0] = 'RUSIZI'
$jan23Data[19,1] = 'MURURU'
$jan23Data[19,2] = 81
$jan23Data[19,3] = 5
$jan23Data[19,4] = 6.17
$jan23Data[19,5] = 'Low'
$jan23Data[20,0] = 'MUSANZE'
$jan23Data[20,1] = 'GASHAKI'
$jan23Data[20,2] = 99
$jan23Data[20,3] = 7
$jan23Data[20,4] = 7.07
$jan23Data[20,5] = 'Low'
$jan23Data[21,0] = 'RULINDO'
$jan23Data[21,1] = 'TUMBA'
$jan23Data[21,2] = 153
$jan23Data[21,3] = 14
$jan23Data[21,4] = 9.15
$jan23Data[21,5] = 'Low'
$jan23Data[22,0] = 'GISAGARA'
$jan23Data[22,1] = 'KIBIRIZI'
$jan23Data[22,2] = 341
$jan23Data[22,3] = 33
$jan23Data[22,4] = 9.68
$jan23Data[22,5] = 'Low'
$jan23Data[23,0] = 'GISAGARA'
$jan23Data[23,1] = 'SAVE'
$jan23Data[23,2] = 251
$jan23Data[23,3] = 25
$jan23Data[23,4] = 9.9600000000000009
$jan23Data[23,5] = 'Low'
$jan23Data[24,0] = 'RUHANGO'
$jan23Data[24,1] = 'BYIMANA'
$jan23Data[24,2] = 70
$jan23Data[24,3] = 7
$jan23Data[24,4] = 10
$jan23Data[24,5] = 'Low'
$jan23Data[25,0] = 'MUHANGA'
$jan23Data[25,1] = 'KABACUZI'
$jan23Data[25,2] = 327
$jan23Data[25,3] = 39
$jan23Data[25,4] = 11.93
$jan23Data[25,5] = 'Low'
$jan23Data[26,0] = 'NYAMASHEKE'
$jan23Data[26,1] = 'BUSHENGE'
$jan23Data[26,2] = 33
$jan23Data[26,3] = 4
$jan23Data[26,4] = 12.12
$jan23Data[26,5] = 'Low'
$jan23Data[27,0] = 'NYAMAGABE'
$jan23Data[27,1] = 'MUSHUBI'
$jan23Data[27,2] = 104
$jan23Data[27,3] = 13
$jan23Data[27,4] = 12.5
$jan23Data[27,5] = 'Low'
$jan23Data[28,0] = 'HUYE'
$jan23Data[28,1] = 'NGOMA'
$jan23Data[28,2] = 23
$jan23Data[28,3] = 3
$jan23Data[28,4] = 13.04
$jan23Data[28,5] = 'Low'
$jan23Data[29,0] = 'NYAMASHEKE'
$jan23Data[29,1] = 'KANJONGO'
$jan23Data[29,2] = 238
$jan23Data[29,3] = 39
$jan23Data[29,4] = 16.39
$jan23Data[29,5] = 'Low'
$jan23Data[30,0] = 'MUSANZE'
$jan23Data[30,1] = 'RWAZA'
$jan23Data[30,2] = 22
$jan23Data[30,3] = 4
$jan23Data[30,4] = 18.18
$jan23Data[30,5] = 'Low'
$jan23Data[31,0] = 'BURERA'
$jan23Data[31,1] = 'CYERU'
$jan23Data[31,2] = 26
$jan23Data[31,3] = 5
$jan23Data[31,4] = 19.23
$jan23Data[31,5] = 'Low'
$jan23Data[32,0] = 'KIREHE'
$jan23Data[32,1] = 'GATORE'
$jan23Data[32,2] = 26
$jan23Data[32,3] = 5
$jan23Data[32,4] = 19.23
$jan23Data[32,5] = 'Low'
$jan23Data[33,0] = 'GISAGARA'
$jan23Data[33,1] = 'KANSI'
$jan23Data[33,2] = 115
$jan23Data[33,3] = 23
$jan23Data[33,4] = 20
$jan23Data[33,5] = 'Medium'
$jan23Data[34,0] = 'HUYE'
$jan23Data[34,1] = 'RUSATIRA'
$jan23Data[34,2] = 98
$jan23Data[34,3] = 21
$jan23Data[34,4] = 21.43
$jan23Data[34,5] = 'Medium'
$jan23Data[35,0] = 'RUSIZI'
$jan23Data[35,1] = 'GASHONGA'
$jan23Data[35,2] = 182
$jan23Data[35,3] = 40
$jan23Data[35,4] = 21.98
$jan23Data[35,5] = 'Medium'
$jan23Data[36,0] = 'RULINDO'
$jan23Data[36,1] = 'CYUNGO'
$jan23Data[36,2] = 22
$jan23Data[36,3] = 5
$jan23Data[36,4] = 22.73
$jan23Data[36,5] = 'Medium'
$jan23Data[37,0] = 'RUBAVU'
$jan23Data[37,1] = 'RUGERERO'
$jan23Data[37,2] = 110
$jan23Data[37,3] = 26
$jan23Data[37,4] = 23.64
$jan23Data[37,5] = 'Medium'
$jan23Data[38,0] = 'RUBAVU'
$jan23Data[38,1] = 'NYAMYUMBA'
$jan23Data[38,2] = 22
$jan23Data[38,3] = 6
$jan23Data[38,4] = 27.27
$jan23Data[38,5] = 'Medium'
$jan23Data[39,0] = 'GISAGARA'
$jan23Data[39,1] = 'NYANZA'
$jan23Data[39,2] = 88
$jan23Data[39,3] = 24
$jan23Data[39,4] = 27.27
$jan23Data[39,5] = 'Medium'
$jan23Data[40,0] = 'HUYE'
$jan23Data[40,1] = 'MBAZI'
$jan23Data[40,2] = 248
$jan23Data[40,3] = 69
$jan23Data[40,4] = 27.82
$jan23Data[40,5] = 'Medium'
$jan23Data[41,0] = 'GAKENKE'
$jan23Data[41,1] = 'RUSASA'
$jan23Data[41,2] = 86
$jan23Data[41,3] = 26
$jan23Data[41,4] = 30.23
$jan23Data[41,5] = 'Medium'
$jan23Data[42,0] = 'RULINDO'
$jan23Data[42,1] = 'KINIHIRA'
$jan23Data[42,2] = 157
$jan23Data[42,3] = 48
$jan23Data[42,4] = 30.57
$jan23Data[42,5] = 'Medium'
$jan23Data[43,0] = 'GISAGARA'
$jan23Data[43,1] = 'MUKINDO'
$jan23Data[43,2] = 73
$jan23Data[43,3] = 24
$jan23Data[43,4] = 32.880000000000003
$jan23Data[43,5] = 'Medium'
$jan23Data[44,0] = 'NYABIHU'
$jan23Data[44,1] = 'RUGERA'
$jan23Data[44,2] = 21
$jan23Data[44,3] = 7
$jan23Data[44,4] = 33.33
$jan23Data[44,5] = 'Medium'
$jan23Data[45,0] = 'MUHANGA'
$jan23Data[45,1] = 'SHYOGWE'
$jan23Data[45,2] = 114
$jan23Data[45,3] = 41
$jan23Data[45,4] = 35.96
$jan23Data[45,5] = 'Medium'
$jan23Data[46,0] = 'GAKENKE'
$jan23Data[46,1] = 'RUSHASHI'
$jan23Data[46,2] = 71
$jan23Data[46,3] = 29
$jan23Data[46,4] = 40.85
$jan23Data[46,5] = 'Medium'
$jan23Data[47,0] = 'RUSIZI'
$jan23Data[47,1] = 'NKANKA'
$jan23Data[47,2] = 251
$jan23Data[47,3] = 103
$jan23Data[47,4] = 41.04
$jan23Data[47,5] = 'Medium'
$jan23Data[48,0] = 'MUSANZE'
$jan23Data[48,1] = 'MUKO'
$jan23Data[48,2] = 24
$jan23Data[48,3] = 10
$jan23Data[48,4] = 41.67
$jan23Data[48,5] = 'Medium'
$jan23Data[49,0] = 'GATSIBO'
$jan23Data[49,1] = 'GITOKI'
$jan23Data[49,2] = 80
$jan23Data[49,3] = 35
$jan23Data[49,4] = 43.75
$jan23Data[49,5] = 'Medium'
$jan23Data[50,0] = 'NYAMASHEKE'
$jan23Data[50,1] = 'KAGANO'
$jan23Data[50,2] = 51
$jan23Data[50,3] = 23
$jan23Data[50,4] = 45.1
$jan23Data[50,5] = 'Medium'
$jan23Data[51,0] = 'NYAMASHEKE'
$jan23Data[51,1] = 'MACUBA'
$jan23Data[51,2] = 37
$jan23Data[51,3] = 19
$jan23Data[51,4] = 51.35
$jan23Data[51,5] = 'High'
$jan23Data[52,0] = 'BURERA'
$jan23Data[52,1] = 'RUGENGABALI'
$jan23Data[52,2] = 88
$jan23Data[52,3] = 46
$jan23Data[52,4] = 52.27
$jan23Data[52,5] = 'High'
$jan23Data[53,0] = 'KAYONZA'
$jan23Data[53,1] = 'MUKARANGE'
$jan23Data[53,2] = 49
$jan23Data[53,3] = 28
$jan23Data[53,4] = 57.14
$jan23Data[53,5] = 'High'
$jan23Data[54,0] = 'RUSIZI'
$jan23Data[54,1] = 'NYAKABUYE'
$jan23Data[54,2] = 23
$jan23Data[54,3] = 14
$jan23Data[54,4] = 60.87
$jan23Data[54,5] = 'High'
$jan23Data[55,0] = 'KAMONYI'
$jan23Data[55,1] = 'MUGINA'
$jan23Data[55,2] = 63
$jan23Data[55,3] = 42
$jan23Data[55,4] = 66.67
$jan23Data[55,5] = 'High'
$jan23Data[56,0] = 'NYAMASHEKE'
$jan23Data[56,1] = 'BUSHEKERI'
$jan23Data[56,2] = 29
$jan23Data[56,3] = 20
$jan23Data[56,4] = 68.97
$jan23Data[56,5] = 'High'
$jan23Data[57,0] = 'NYAMAGABE'
$jan23Data[57,1] = 'MBAZI'
$jan23Data[57,2] = 42
$jan23Data[57,3] = 30
$jan23Data[57,4] = 71.430000000000007
$jan23Data[57,5] = 'High'
$jan23Data[58,0] = 'GICUMBI'
$jan23Data[58,1] = 'BWISIGE'
$jan23Data[58,2] = 36
$jan23Data[58,3] = 26
$jan23Data[58,4] = 72.22
$jan23Data[58,5] = 'High'
$jan23Data[59,0] = 'MUHANGA'
$jan23Data[59,1] = 'KIBANGU'
$jan23Data[59,2] = 96
$jan23Data[59,3] = 71
$jan23Data[59,4] = 73.959999999999994
$jan23Data[59,5] = 'High'
$jan23Data[60,0] = 'KAYONZA'
$jan23Data[60,1] = 'KABARONDO'
$jan23Data[60,2] = 62
$jan23Data[60,3] = 46
$jan23Data[60,4] = 74.19
$jan23Data[60,5] = 'High'
$jan23Data[61,0] = 'RUTSIRO'
$jan23Data[61,1] = 'KIGEYO'
$jan23Data[61,2] = 47
$jan23Data[61,3] = 36
$jan23Data[61,4] = 76.599999999999994
$jan23Data[61,5] = 'High'
$jan23Data[62,0] = 'RUTSIRO'
$jan23Data[62,1] = 'MUSHUBATI'
$jan23Data[62,2] = 46
$jan23Data[62,3] = 42
$jan23Data[62,4] = 91.3
$jan23Data[62,5] = 'High'
$wsJan23.Range("A1:F63").Value = $jan23Data

# ---- Jan_23 sheet: column widths ----
$wsJan23.Columns.Item(1).ColumnWidth = 14.833333333333334
$wsJan23.Columns.Item(2).ColumnWidth = 12.0
$wsJan23.Columns.Item(3).ColumnWidth = 12.333333333333334
$wsJan23.Columns.Item(4).ColumnWidth = 12.0
$wsJan23.Columns.Item(5).ColumnWidth = 13.666666666666666

# ---- Selections per sheet (also drives tabSelected via Activate order) ----
$wsJune.Activate()
$wsJune.Range("G1").Select()

$wsSept.Activate()
$wsSept.Range("I22").Select()

$wsMar23.Activate()
$wsMar23.Range("D27").Select()

$wsJune23.Activate()
$wsJune23.Range("A1").Select()

$wsSept23.Activate()
$wsSept23.Range("G27").Select()

$wsJan24.Activate()
$wsJan24.Range("I27").Select()

# Jan_23 is the active/selected sheet at the end (matches activeTab=2, tabSelected on sheet3)
$wsJan23.Activate()
$wsJan23.Range("E1:E1048576").Select()

